# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal").
# Insert 3 new rows of data (the latest week's prices) right before the
# existing row 786, shifting all subsequent rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before current row 786 (pushes old rows 786+ down to 789+)
$ws.Rows("786:788").Insert()

# Row 786 - 1a amarillo
$ws.Cells.Item(786, 4).Value = 45077
$ws.Cells.Item(786, 12).Value = "1a amarillo"
$ws.Cells.Item(786, 13).Value = 1050
$ws.Cells.Item(786, 14).Value = 5800
$ws.Cells.Item(786, 15).Value = 6000
$ws.Cells.Item(786, 16).Value = 5900
$ws.Cells.Item(786, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(786, 19).Value = 328
$ws.Cells.Item(786, 20).Value = 18

# Row 787 - 2a amarillo
$ws.Cells.Item(787, 4).Value = 45077
$ws.Cells.Item(787, 12).Value = "2a amarillo"
$ws.Cells.Item(787, 13).Value = 780
$ws.Cells.Item(787, 14).Value = 3800
$ws.Cells.Item(787, 15).Value = 4000
$ws.Cells.Item(787, 16).Value = 3900
$ws.Cells.Item(787, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(787, 19).Value = 217
$ws.Cells.Item(787, 20).Value = 18

# Row 788 - 3a amarillo
$ws.Cells.Item(788, 4).Value = 45077
$ws.Cells.Item(788, 12).Value = "3a amarillo"
$ws.Cells.Item(788, 13).Value = 450
$ws.Cells.Item(788, 14).Value = 1800
$ws.Cells.Item(788, 15).Value = 2000
$ws.Cells.Item(788, 16).Value = 1900
$ws.Cells.Item(788, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(788, 19).Value = 106
$ws.Cells.Item(788, 20).Value = 18

# Fill in remaining columns for the 3 new rows by copying from the row
# that is now directly below each new row's corresponding "family" (they
# share the same market/product/origin metadata as the rest of the sheet).
$cols = @(1,2,3,5,6,7,8,9,10,11,18)
foreach ($r in 786..788) {
    $src = $r + 3
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($src, $c).Value2
    }
}
